$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5535.8096
$ws.Range("I18").Value = 2658.2666
$ws.Range("J18").Value = 12729.667
$ws.Range("K18").Value = 2658.2666
$ws.Range("L18").Value = 12729.667
$ws.Range("M18").Value = -2374.2666
$ws.Range("N18").Value = -13297.667
$ws.Range("H40").Value = 1965
$ws.Range("I40").Value = 2263.3333
$ws.Range("J40").Value = 1666.6666
$ws.Range("K40").Value = 2263.3333
$ws.Range("L40").Value = 1666.6666
$ws.Range("M40").Value = -2088.3333
$ws.Range("N40").Value = -2016.6666
$ws.Range("H74").Value = 2999.5
$ws.Range("I74").Value = 2999.5
$ws.Range("K74").Value = 2999.5
$ws.Range("M74").Value = -2063.5
$ws.Range("H77").Value = 2999.5
$ws.Range("I77").Value = 2999.5
$ws.Range("K77").Value = 14997.5
$ws.Range("M77").Value = -10317.5
$ws.Range("H98").Value = 944.561
$ws.Range("I98").Value = 762.94116
$ws.Range("K98").Value = 762.94116
$ws.Range("M98").Value = 735.05884
$ws.Range("H116").Value = 13949.5
$ws.Range("I116").Value = 34499.668
$ws.Range("K116").Value = 34499.668
$ws.Range("M116").Value = -31057.668
$ws.Range("H122").Value = 944.561
$ws.Range("I122").Value = 762.94116
$ws.Range("K122").Value = 2288.82348
$ws.Range("M122").Value = 161.17652
$ws.Range("H127").Value = 2741.3635
$ws.Range("I127").Value = 2706.875
$ws.Range("J127").Value = 2833.3333
$ws.Range("K127").Value = 8120.625
$ws.Range("L127").Value = 8499.999899999999
$ws.Range("M127").Value = -3160.625
$ws.Range("N127").Value = -18419.9999
$ws.Range("H132").Value = 1357.2106
$ws.Range("I132").Value = 1253.6
$ws.Range("K132").Value = 3760.8
$ws.Range("M132").Value = -1230.8
$ws.Range("H138").Value = 3860.6155
$ws.Range("I138").Value = 7929
$ws.Range("J138").Value = 2640.1
$ws.Range("K138").Value = 23787
$ws.Range("L138").Value = 7920.299999999999
$ws.Range("M138").Value = -18647
$ws.Range("N138").Value = -18200.3
$ws.Range("H140").Value = 64272.473
$ws.Range("J140").Value = 64272.473
$ws.Range("L140").Value = 64272.473
$ws.Range("N140").Value = -74632.473

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5042.784
$ws.Range("I32").Value = 3897.4048
$ws.Range("J32").Value = 10387.889
$ws.Range("K32").Value = 3897.4048
$ws.Range("L32").Value = 10387.889
$ws.Range("M32").Value = -3610.4048
$ws.Range("N32").Value = -10961.889
$ws.Range("H132").Value = 1486.4333
$ws.Range("I132").Value = 1239.4546
$ws.Range("J132").Value = 2165.625
$ws.Range("K132").Value = 3718.3638
$ws.Range("L132").Value = 6496.875
$ws.Range("M132").Value = -1188.3638
$ws.Range("N132").Value = -11556.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3491.375
$ws.Range("I20").Value = 3733
$ws.Range("J20").Value = 3249.75
$ws.Range("K20").Value = 3733
$ws.Range("L20").Value = 3249.75
$ws.Range("M20").Value = -3486
$ws.Range("N20").Value = -3743.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2753.889
$ws.Range("I31").Value = 2633.111
$ws.Range("J31").Value = 2874.6667
$ws.Range("K31").Value = 2633.111
$ws.Range("L31").Value = 2874.6667
$ws.Range("M31").Value = -2338.111
$ws.Range("N31").Value = -3464.6667
$ws.Range("H34").Value = 2753.889
$ws.Range("I34").Value = 2633.111
$ws.Range("J34").Value = 2874.6667
$ws.Range("K34").Value = 2633.111
$ws.Range("L34").Value = 2874.6667
$ws.Range("M34").Value = -2431.111
$ws.Range("N34").Value = -3278.6667
$ws.Range("H58").Value = 3662.4666
$ws.Range("J58").Value = 4665.5557
$ws.Range("L58").Value = 4665.5557
$ws.Range("N58").Value = -5071.5557
$ws.Range("H99").Value = 3413.1428
$ws.Range("I99").Value = 2223
$ws.Range("K99").Value = 2223
$ws.Range("M99").Value = -725
$ws.Range("H126").Value = 3413.1428
$ws.Range("I126").Value = 2223
$ws.Range("K126").Value = 6669
$ws.Range("M126").Value = -4199
$ws.Range("H132").Value = 3391.6316
$ws.Range("I132").Value = 2637
$ws.Range("J132").Value = 4429.25
$ws.Range("K132").Value = 7911
$ws.Range("L132").Value = 13287.75
$ws.Range("M132").Value = -5381
$ws.Range("N132").Value = -18347.75
$ws.Range("H136").Value = 3662.4666
$ws.Range("J136").Value = 4665.5557
$ws.Range("L136").Value = 13996.6671
$ws.Range("N136").Value = -19096.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 182.3077
$ws.Range("I2").Value = 144.66667
$ws.Range("J2").Value = 267
$ws.Range("K2").Value = 868.0000200000001
$ws.Range("L2").Value = 1602
$ws.Range("M2").Value = -755.0000200000001
$ws.Range("N2").Value = -1828
$ws.Range("H5").Value = 386.27777
$ws.Range("I5").Value = 314.63635
$ws.Range("J5").Value = 498.85715
$ws.Range("K5").Value = 943.90905
$ws.Range("L5").Value = 1496.57145
$ws.Range("M5").Value = -831.90905
$ws.Range("N5").Value = -1720.57145
$ws.Range("H70").Value = 3502.4
$ws.Range("J70").Value = 5333.3335
$ws.Range("L70").Value = 16000.0005
$ws.Range("N70").Value = -16630.0005
$ws.Range("H73").Value = 3502.4
$ws.Range("J73").Value = 5333.3335
$ws.Range("L73").Value = 16000.0005
$ws.Range("N73").Value = -18184.0005
$ws.Range("H107").Value = 454.6111
$ws.Range("J107").Value = 452.2
$ws.Range("L107").Value = 1356.6
$ws.Range("N107").Value = -5196.6
$ws.Range("H131").Value = 26579.777
$ws.Range("J131").Value = 29807.25
$ws.Range("L131").Value = 89421.75
$ws.Range("N131").Value = -99501.75
$ws.Range("H135").Value = 386.27777
$ws.Range("I135").Value = 314.63635
$ws.Range("J135").Value = 498.85715
$ws.Range("K135").Value = 2831.72715
$ws.Range("L135").Value = 4489.71435
$ws.Range("M135").Value = -296.7271499999997
$ws.Range("N135").Value = -9559.71435
$ws.Range("H140").Value = 2546.6
$ws.Range("I140").Value = 1344.6875
$ws.Range("J140").Value = 3920.2144
$ws.Range("K140").Value = 4034.0625
$ws.Range("L140").Value = 11760.6432
$ws.Range("M140").Value = 1145.9375
$ws.Range("N140").Value = -22120.6432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3820.3572
$ws.Range("I102").Value = 4546.353
$ws.Range("J102").Value = 2698.3635
$ws.Range("K102").Value = 4546.353
$ws.Range("L102").Value = 2698.3635
$ws.Range("M102").Value = -2924.353
$ws.Range("N102").Value = -5942.363499999999
$ws.Range("H126").Value = 33039.848
$ws.Range("I126").Value = 3339.2
$ws.Range("K126").Value = 10017.6
$ws.Range("M126").Value = -7547.599999999999
$ws.Range("H132").Value = 3560.9565
$ws.Range("I132").Value = 3015.923
$ws.Range("J132").Value = 4269.5
$ws.Range("K132").Value = 9047.769
$ws.Range("L132").Value = 12808.5
$ws.Range("M132").Value = -6517.769
$ws.Range("N132").Value = -17868.5
$ws.Range("H140").Value = 47511
$ws.Range("J140").Value = 47511
$ws.Range("L140").Value = 47511
$ws.Range("N140").Value = -57871

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9377
$ws.Range("I16").Value = 9377
$ws.Range("K16").Value = 9377
$ws.Range("M16").Value = -9207
$ws.Range("H40").Value = 15958.111
$ws.Range("J40").Value = 8689.666999999999
$ws.Range("L40").Value = 8689.666999999999
$ws.Range("N40").Value = -8961.666999999999
$ws.Range("H61").Value = 2462.9092
$ws.Range("I61").Value = 1761.5
$ws.Range("J61").Value = 4333.3335
$ws.Range("K61").Value = 1761.5
$ws.Range("L61").Value = 4333.3335
$ws.Range("M61").Value = -1559.5
$ws.Range("N61").Value = -4737.3335
$ws.Range("H113").Value = 2462.9092
$ws.Range("I113").Value = 1761.5
$ws.Range("J113").Value = 4333.3335
$ws.Range("K113").Value = 1761.5
$ws.Range("L113").Value = 4333.3335
$ws.Range("M113").Value = 408.5
$ws.Range("N113").Value = -8673.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 455.45456
$ws.Range("I113").Value = 344.83334
$ws.Range("J113").Value = 588.2
$ws.Range("K113").Value = 1034.50002
$ws.Range("L113").Value = 1764.6
$ws.Range("M113").Value = 1135.49998
$ws.Range("N113").Value = -6104.6
$ws.Range("H123").Value = 47635.2
$ws.Range("J123").Value = 47635.2
$ws.Range("L123").Value = 47635.2
$ws.Range("N123").Value = -57435.2
$ws.Range("H126").Value = 10825.083
$ws.Range("I126").Value = 10900.091
$ws.Range("M126").Value = -30230.273
$ws.Range("H132").Value = 1777.5
$ws.Range("I132").Value = 1065.619
$ws.Range("J132").Value = 3136.5454
$ws.Range("K132").Value = 3196.857
$ws.Range("L132").Value = 9409.636200000001
$ws.Range("M132").Value = -666.857
$ws.Range("N132").Value = -14469.6362
